$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# Header row (row 1): 工号 / 金额 / 名称 / 发放时间  (C1 changes 名称 -> name col)
# ---------------------------------------------------------------------------
$ws.Range("A1").Value = "工号"
$ws.Range("B1").Value = "金额"
$ws.Range("C1").Value = "名称"
$ws.Range("D1").Value = "发放时间"

# ---------------------------------------------------------------------------
# Data rows 2-5: one employee (130259), amount 100, four benefit kinds,
# all issued 2020-09
# ---------------------------------------------------------------------------
$ws.Range("A2").Value = "130259"
$ws.Range("B2").Value = "100"
$ws.Range("C2").Value = "半年奖"
$ws.Range("D2").Value = "2020-09"

$ws.Range("A3").Value = "130259"
$ws.Range("B3").Value = "100"
$ws.Range("C3").Value = "年终奖"
$ws.Range("D3").Value = "2020-09"

$ws.Range("A4").Value = "130259"
$ws.Range("B4").Value = "100"
$ws.Range("C4").Value = "十三薪"
$ws.Range("D4").Value = "2020-09"

$ws.Range("A5").Value = "130259"
$ws.Range("B5").Value = "100"
$ws.Range("C5").Value = "冷饮费"
$ws.Range("D5").Value = "2020-09"

# Columns A:C only carry a column-level default style (no alignment); force
# the same "text, centered" look the header/first data row already uses so
# the newly written rows 3-5 line up with rows 1-2.
$ws.Range("A2:C5").NumberFormat = "@"
$ws.Range("A2:C5").HorizontalAlignment = -4108

# ---------------------------------------------------------------------------
# Remark note in column F, row 1 - red text, left aligned, wide column
# ---------------------------------------------------------------------------
$ws.Range("F1").Value = "备注：半年奖+年终奖=奖金；除了奖金，十三薪，其他的都是福利；"
$ws.Range("F1").Font.Color = 255
$ws.Range("F1").NumberFormat = "@"
$ws.Range("F1").HorizontalAlignment = -4131
$ws.Columns("F").ColumnWidth = 75.375

# ---------------------------------------------------------------------------
# Selection marker as left by the author when they saved the file
# ---------------------------------------------------------------------------
$ws.Range("F10").Select() | Out-Null
